{"js": "// The document contains a 5-column table of \"two-digit \u00f7 one-digit\"\n// division prompts (e.g. \"93\u00f77=\"). The commit replaces the 25 filled-in\n// prompts with a new set of values, in the same left-to-right,\n// top-to-bottom reading order. Empty placeholder cells are untouched.\nconst newValues = [\n  \"38\u00f78=\", \"42\u00f79=\", \"77\u00f74=\", \"94\u00f77=\", \"18\u00f76=\",\n  \"99\u00f75=\", \"69\u00f79=\", \"18\u00f76=\", \"59\u00f74=\", \"80\u00f79=\",\n  \"39\u00f78=\", \"40\u00f74=\", \"66\u00f76=\", \"76\u00f79=\", \"75\u00f75=\",\n  \"80\u00f76=\", \"74\u00f79=\", \"77\u00f73=\", \"25\u00f72=\", \"75\u00f76=\",\n  \"97\u00f73=\", \"98\u00f74=\", \"34\u00f72=\", \"11\u00f72=\", \"40\u00f78=\"\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (const p of paragraphs.items) {\n  p.load(\"text\");\n}\nawait context.sync();\n\n// Walk the document's paragraphs in order, picking out the ones that hold\n// a division prompt (\"digits \u00f7 digits =\"); replace their text in place\n// (via the paragraph's own range) so run/paragraph formatting (fonts,\n// size, justification, ...) is preserved untouched.\nconst exprPattern = /^\\s*\\d+\u00f7\\d+=\\s*$/;\nlet idx = 0;\nfor (const p of paragraphs.items) {\n  if (idx >= newValues.length) break;\n  if (exprPattern.test(p.text)) {\n    const range = p.getRange();\n    range.insertText(newValues[idx], Word.InsertLocation.replace);\n    idx++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document contains a 5-column table of \"two-digit \u00f7 one-digit\"\n# division prompts (e.g. \"93\u00f77=\"). The commit replaces the 25 filled-in\n# prompts with a new set of values, in the same left-to-right,\n# top-to-bottom reading order. Empty placeholder cells are left alone.\n$newValues = @(\n  \"38\u00f78=\", \"42\u00f79=\", \"77\u00f74=\", \"94\u00f77=\", \"18\u00f76=\",\n  \"99\u00f75=\", \"69\u00f79=\", \"18\u00f76=\", \"59\u00f74=\", \"80\u00f79=\",\n  \"39\u00f78=\", \"40\u00f74=\", \"66\u00f76=\", \"76\u00f79=\", \"75\u00f75=\",\n  \"80\u00f76=\", \"74\u00f79=\", \"77\u00f73=\", \"25\u00f72=\", \"75\u00f76=\",\n  \"97\u00f73=\", \"98\u00f74=\", \"34\u00f72=\", \"11\u00f72=\", \"40\u00f78=\"\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Walk every cell of the table in document order (row-major), picking out\n# the ones that hold a division prompt (\"digits \u00f7 digits =\"); replace the\n# cell's text in place so paragraph/run formatting (fonts, size,\n# justification, ...) is preserved untouched. A cell's .Range.Text carries\n# trailing end-of-cell/paragraph marks (chr 13 / chr 7), so trim those\n# before testing the pattern.\n$idx = 0\nforeach ($cell in $t.Range.Cells) {\n  if ($idx -ge $newValues.Count) { break }\n  $txt = $cell.Range.Text.TrimEnd([char]13, [char]7)\n  if ($txt -match '^\\s*\\d+\u00f7\\d+=\\s*$') {\n    $cell.Range.Text = $newValues[$idx]\n    $idx++\n  }\n}\n"}
